# Add the new monthly data row (row 65) to the national CPI sheet:
# interannual + monthly variation figures for the latest period (2022-04-01).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 65
$values = @(6, 5.9, 3.3, 9.9, 4.6, 5.5, 6.4, 5.3, 3.7, 5.2, 3.7, 7.3, 5.3)

for ($i = 0; $i -lt $values.Length; $i++) {
    $col = $i + 1
    $ws.Cells.Item($row, $col).Value = $values[$i]
}

# Column N holds the period date (serial 44652 = 2022-04-01); it already
# inherits the yyyy-mm-dd date style from the column definition.
$ws.Cells.Item($row, 14).Value = 44652
